$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 333.6154
$ws.Range("I2").Value = 318
$ws.Range("J2").Value = 368.75
$ws.Range("K2").Value = 318
$ws.Range("L2").Value = 368.75
$ws.Range("M2").Value = -205
$ws.Range("N2").Value = -594.75

$ws.Range("H40").Value = 6423.75
$ws.Range("I40").Value = 3597.5
$ws.Range("K40").Value = 3597.5
$ws.Range("M40").Value = -3422.5

$ws.Range("H42").Value = 220.8
$ws.Range("J42").Value = 524.5
$ws.Range("L42").Value = 1573.5
$ws.Range("N42").Value = -2033.5

$ws.Range("H55").Value = 220.85715
$ws.Range("I55").Value = 97.666664
$ws.Range("K55").Value = 97.666664
$ws.Range("M55").Value = 116.333336

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 696.7143
$ws.Range("I25").Value = 762.8333
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 762.8333
$ws.Range("L25").Value = 300
$ws.Range("M25").Value = -360.8333
$ws.Range("N25").Value = -1104

$ws.Range("H32").Value = 2944.6
$ws.Range("I32").Value = 3074.652
$ws.Range("K32").Value = 3074.652
$ws.Range("M32").Value = -2787.652

$ws.Range("H132").Value = 4151.077
$ws.Range("I132").Value = 3346.4
$ws.Range("K132").Value = 10039.2
$ws.Range("M132").Value = -7509.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 10032.8
$ws.Range("I54").Value = 10032.8
$ws.Range("K54").Value = 10032.8
$ws.Range("M54").Value = -9548.799999999999

$ws.Range("H82").Value = 15522.6
$ws.Range("I82").Value = 15522.6
$ws.Range("K82").Value = 15522.6
$ws.Range("M82").Value = -15139.6

$ws.Range("H85").Value = 15522.6
$ws.Range("I85").Value = 15522.6
$ws.Range("K85").Value = 15522.6
$ws.Range("M85").Value = -14196.6

$ws.Range("H86").Value = 1749.75
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -4245.5

$ws.Range("H89").Value = 1749.75
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -21229.5

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 999
$ws.Range("I7").Value = 999
$ws.Range("K7").Value = 999
$ws.Range("M7").Value = -886

$ws.Range("H35").Value = 166.5
$ws.Range("I35").Value = 166.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 166.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 127.5
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H7").Value = 478
$ws.Range("I7").Value = 97.5
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 292.5
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -180.5
$ws.Range("N7").Value = -6224

$ws.Range("H131").Value = 1674
$ws.Range("I131").Value = 1157.8572
$ws.Range("K131").Value = 3473.5716
$ws.Range("M131").Value = 1566.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19771
$ws.Range("J57").Value = 23450
$ws.Range("L57").Value = 23450
$ws.Range("N57").Value = -25090

$ws.Range("H70").Value = 4749.5
$ws.Range("I70").Value = 4499.5
$ws.Range("K70").Value = 4499.5
$ws.Range("M70").Value = -4229.5

$ws.Range("H73").Value = 4749.5
$ws.Range("I73").Value = 4499.5
$ws.Range("K73").Value = 4499.5
$ws.Range("M73").Value = -3563.5

$ws.Range("H122").Value = 1052
$ws.Range("J122").Value = 1424.5
$ws.Range("L122").Value = 4273.5
$ws.Range("N122").Value = -9173.5

$ws.Range("H126").Value = 3367.5715
$ws.Range("I126").Value = 3899.75
$ws.Range("J126").Value = 2658
$ws.Range("K126").Value = 11699.25
$ws.Range("L126").Value = 7974
$ws.Range("M126").Value = -9229.25
$ws.Range("N126").Value = -12914

$ws.Range("H132").Value = 5380.9
$ws.Range("J132").Value = 8999.5
$ws.Range("L132").Value = 26998.5
$ws.Range("N132").Value = -32058.5

$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H138").Value = 99429
$ws.Range("J138").Value = 99429
$ws.Range("L138").Value = 99429
$ws.Range("N138").Value = -109709

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1795.3334
$ws.Range("I126").Value = 1795.3334
$ws.Range("K126").Value = 5386.0002
$ws.Range("M126").Value = -2916.0002

$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -105120
